$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Current (before) column layout on Sheet1:
#   A query_id | B gold_standard | C GPT-5_instant | D GPT-4o |
#   E MetaAI_Llama-4 | F Gemini_2.5-Flash | G Claude_Sonnet-4.5 | H Mistral_7B |
#   I GPT-5_thinking | J Deepseek | K Deepseek_DeepThink | L Deepseek_error
#
# Target (after) column layout:
#   A query_id | B gold_standard | C GPT-5_instant | D GPT-4o |
#   E Deepseek | F Gemini_2.5-Flash | G Claude_Sonnet-4.5 | H Mistral_7B
#
# i.e. the "MetaAI_Llama-4", "GPT-5_thinking", "Deepseek_DeepThink" and
# "Deepseek_error" columns (and their answer columns) are removed, and the
# "Deepseek" column is moved in to sit right after "GPT-4o".
# ---------------------------------------------------------------------------

# 1) Capture the "Deepseek" column (J) contents before anything is deleted,
#    since it needs to survive and move left into column E.
$deepseekHeader = $ws.Range("J1").Value()
$deepseek2 = $ws.Range("J2").Value()
$deepseek3 = $ws.Range("J3").Value()
$deepseek4 = $ws.Range("J4").Value()
$deepseek5 = $ws.Range("J5").Value()

# 2) Drop the trailing block of columns that are no longer needed in one
#    single delete: I (GPT-5_thinking), J (Deepseek - already saved above),
#    K (Deepseek_DeepThink) and L (Deepseek_error).
$ws.Columns("I:L").Delete()

# 3) Column E currently still holds "MetaAI_Llama-4" -- overwrite it in
#    place with the saved "Deepseek" column data (keeps the existing
#    formatting/style of column E, e.g. the wrap-text style on the data
#    rows, while swapping the actual content).
$ws.Range("E1").Value = $deepseekHeader
$ws.Range("E2").Value = $deepseek2
$ws.Range("E3").Value = $deepseek3
$ws.Range("E4").Value = $deepseek4
$ws.Range("E5").Value = $deepseek5

# 4) Restore the selection to match the post-edit workbook state.
$ws.Range("G2").Select()
